$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace D2/E2 values with the new special tags, drop F2 entirely
$ws.Range("D2").Value = '${EMPTY}'
$ws.Range("E2").Value = '${NULL}'
$ws.Range("F2").ClearContents()

# Update the active selection on the sheet view
$ws.Range("E12").Select()
